# Auto-generated script: apply 'Generate Report for Handback' edit
$wb = $excel.ActiveWorkbook

# Key filename / URL strings
$da86563b_md = "da86563b-6ed4-43e2-99b3-838d7a2616c7ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$e2e_da86563b_md = "e2e\da86563b-6ed4-43e2-99b3-838d7a2616c7ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$b399af61_md = "b399af61-778a-42f6-b851-e9c57ab1dc14ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$e2e_b399af61_md = "e2e\b399af61-778a-42f6-b851-e9c57ab1dc14ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$da86563b_zhcn_xlf = "da86563b-6ed4-43e2-99b3-838d7a2616c7oooooooooooooooooooooooooooooooooooooooo.86cd9389a8ac492b4524e8338659dae39e038732.zh-cn.xlf"
$b399af61_zhcn_xlf = "b399af61-778a-42f6-b851-e9c57ab1dc14oooooooooooooooooooooooooooooooooooooooo.e8ca229ebfcc611948cb21c69897002b596676ac.zh-cn.xlf"
$da86563b_dede_xlf = "da86563b-6ed4-43e2-99b3-838d7a2616c7oooooooooooooooooooooooooooooooooooooooo.86cd9389a8ac492b4524e8338659dae39e038732.de-de.xlf"
$b399af61_dede_xlf = "b399af61-778a-42f6-b851-e9c57ab1dc14oooooooooooooooooooooooooooooooooooooooo.e8ca229ebfcc611948cb21c69897002b596676ac.de-de.xlf"
$da86563b_url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58d3bff3135039288a000dc49dae2b6458488f86/e2e/da86563b-6ed4-43e2-99b3-838d7a2616c7ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$b399af61_url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2333acaa1d7dc7ebaf3f112e8dc56a1a3248604/e2e/b399af61-778a-42f6-b851-e9c57ab1dc14ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"

# ---- Sheet: Overview ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()
$ws.Range("A2").Value = $b399af61_md
$ws.Range("B2").Value = $e2e_b399af61_md
$ws.Range("C2").Value = "'.md"
$ws.Range("E2").Value = "'Handed back: in sync with en-US"
$ws.Range("F2").Value = "'Handed back: in sync with en-US"
$ws.Range("G2").Value = "'2016-08-21 00:36:35"
$ws.Range("A3").Value = $da86563b_md
$ws.Range("B3").Value = $e2e_da86563b_md
$ws.Range("C3").Value = "'.md"
$ws.Range("E3").Value = "'In Translation"
$ws.Range("F3").Value = "'In Translation"
$ws.Range("G3").Value = "'2016-08-21 00:36:03"

# ---- Sheet: zh-cn ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()
$ws.Range("A2").Value = $b399af61_md
$ws.Range("B2").Value = "'.md"
$ws.Range("C2").Value = "'Handed back: in sync with en-US"
$ws.Range("D2").Value = "'e2e"
$ws.Range("E2").Value = "'ht"
$ws.Range("F2").Value = "'False"
$ws.Range("G2").Value = $b399af61_zhcn_xlf
$ws.Range("H2").Value = "'2016-08-21 00:36:31"
$ws.Range("I2").Value = $b399af61_md
$ws.Range("J2").Value = $b399af61_zhcn_xlf
$ws.Range("K2").Value = "'2016-08-21 00:36:48"
$ws.Range("M2").Value = "'True"
$ws.Range("O2").Value = "'False"
$ws.Range("A3").Value = $da86563b_md
$ws.Range("B3").Value = "'.md"
$ws.Range("C3").Value = "'In Translation"
$ws.Range("D3").Value = "'e2e"
$ws.Range("E3").Value = "'ht"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = $da86563b_zhcn_xlf
$ws.Range("H3").Value = "'2016-08-21 00:35:56"
$ws.Range("K3").Value = "'0001-01-01 00:00:00"
$ws.Range("M3").Value = "'True"
$ws.Range("O3").Value = "'False"

# ---- Sheet: de-de ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()
$ws.Range("A2").Value = $b399af61_md
$ws.Range("B2").Value = "'.md"
$ws.Range("C2").Value = "'Handed back: in sync with en-US"
$ws.Range("D2").Value = "'e2e"
$ws.Range("E2").Value = "'ht"
$ws.Range("F2").Value = "'False"
$ws.Range("G2").Value = $b399af61_dede_xlf
$ws.Range("H2").Value = "'2016-08-21 00:36:35"
$ws.Range("I2").Value = $b399af61_md
$ws.Range("J2").Value = $b399af61_dede_xlf
$ws.Range("K2").Value = "'2016-08-21 00:36:54"
$ws.Range("M2").Value = "'True"
$ws.Range("O2").Value = "'False"
$ws.Range("A3").Value = $da86563b_md
$ws.Range("B3").Value = "'.md"
$ws.Range("C3").Value = "'In Translation"
$ws.Range("D3").Value = "'e2e"
$ws.Range("E3").Value = "'ht"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = $da86563b_dede_xlf
$ws.Range("H3").Value = "'2016-08-21 00:36:03"
$ws.Range("K3").Value = "'0001-01-01 00:00:00"
$ws.Range("M3").Value = "'True"
$ws.Range("O3").Value = "'False"


# ---- Hyperlinks ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Add($ws.Range("B2"), $da86563b_url, "", "", $e2e_b399af61_md) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $b399af61_url, "", "", $e2e_da86563b_md) | Out-Null
$ws.Range("B2").Font.Underline = 1
$ws.Range("B2").Font.Color = 15570276
$ws.Range("B3").Font.Underline = 1
$ws.Range("B3").Font.Color = 15570276

$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Add($ws.Range("A2"), $da86563b_url, "", "", $b399af61_md) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), $b399af61_url, "", "", $b399af61_md) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $da86563b_url, "", "", $da86563b_md) | Out-Null
$ws.Range("A2").Font.Underline = 1
$ws.Range("A2").Font.Color = 15570276
$ws.Range("I2").Font.Underline = 1
$ws.Range("I2").Font.Color = 15570276
$ws.Range("A3").Font.Underline = 1
$ws.Range("A3").Font.Color = 15570276

$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Add($ws.Range("A2"), $da86563b_url, "", "", $b399af61_md) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), $b399af61_url, "", "", $b399af61_md) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $da86563b_url, "", "", $da86563b_md) | Out-Null
$ws.Range("A2").Font.Underline = 1
$ws.Range("A2").Font.Color = 15570276
$ws.Range("I2").Font.Underline = 1
$ws.Range("I2").Font.Color = 15570276
$ws.Range("A3").Font.Underline = 1
$ws.Range("A3").Font.Color = 15570276

# ---- Column widths ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Columns.Item(5).ColumnWidth = 29.15
$ws.Columns.Item(6).ColumnWidth = 29.15

$ws = $wb.Worksheets.Item("zh-cn")
$ws.Columns.Item(3).ColumnWidth = 29.15
$ws.Columns.Item(9).ColumnWidth = 39.15
$ws.Columns.Item(10).ColumnWidth = 39.15

$ws = $wb.Worksheets.Item("de-de")
$ws.Columns.Item(3).ColumnWidth = 29.15
$ws.Columns.Item(9).ColumnWidth = 39.15
$ws.Columns.Item(10).ColumnWidth = 39.15

